# SSD Team Skills Inventory - fill in the previously-empty "requirements analysis"
# (column D) ratings, bump a handful of the "No of 1s" tallies in column H, and
# leave the selection on the last cell that was touched (H11), matching the author's
# worked session.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 was hand-formatted slightly differently from the rest of the column (green fill,
# centered horizontally, but no border/vertical-centering) -- reproduce that exactly by
# clearing the border + vertical alignment it would otherwise inherit.
$d3 = $ws.Range("D3")
$d3.Borders.LineStyle = -4142
$d3.Interior.Color = 5296274
$d3.HorizontalAlignment = -4108
$d3.VerticalAlignment = -4107
$d3.Value = 1

# D4:D18 follow the sheet's existing colour convention used on every other skill column
# (green = 1, white = 2, orange = 3), bordered + centered like their neighbours.
$cell = $ws.Range("D4")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 49407
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 3

$cell = $ws.Range("D5")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 5296274
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 1

$cell = $ws.Range("D6")
$cell.Borders.LineStyle = 1
$cell.Interior.Pattern = -4142
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 2

$cell = $ws.Range("D7")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 5296274
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 1

$cell = $ws.Range("D8")
$cell.Borders.LineStyle = 1
$cell.Interior.Pattern = -4142
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 2

$cell = $ws.Range("D9")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 5296274
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 1

$cell = $ws.Range("D10")
$cell.Borders.LineStyle = 1
$cell.Interior.Pattern = -4142
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 2

$cell = $ws.Range("D11")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 5296274
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 1

$cell = $ws.Range("D12")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 49407
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 3

$cell = $ws.Range("D13")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 49407
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 3

$cell = $ws.Range("D14")
$cell.Borders.LineStyle = 1
$cell.Interior.Pattern = -4142
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 2

$cell = $ws.Range("D15")
$cell.Borders.LineStyle = 1
$cell.Interior.Pattern = -4142
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 2

$cell = $ws.Range("D16")
$cell.Borders.LineStyle = 1
$cell.Interior.Pattern = -4142
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 2

$cell = $ws.Range("D17")
$cell.Borders.LineStyle = 1
$cell.Interior.Color = 5296274
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 1

$cell = $ws.Range("D18")
$cell.Borders.LineStyle = 1
$cell.Interior.Pattern = -4142
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.Value = 2

# Column H ("No of 1s") tallies shift now that column D is no longer blank.
$ws.Range("H3").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("H7").Value = 3
$ws.Range("H9").Value = 3
$ws.Range("H11").Value = 3
$ws.Range("H17").Value = 1

# Leave the cursor where the author last left it.
$ws.Range("H11").Select()

Write-Output "SSD Team Skills Inventory: column D filled in, H tallies refreshed"
